$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Results appended by the R script run: two new price rows (124-125) added
# to the bottom of the BWZ.MI table, following the existing layout:
# date | volume | high | low | open | close | adj_close(text) | ticker

# --- Row 124 ---
# Column A keeps the same date/time number format as the rest of the column.
$ws.Range("A123").Copy()
$ws.Range("A124").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A124").Value = 45499.2916666667

$ws.Range("B124").Value = 0
$ws.Range("C124").Value = 0.680000007152557
$ws.Range("D124").Value = 0.680000007152557
$ws.Range("E124").Value = 0.680000007152557
$ws.Range("F124").Value = 0.680000007152557

# adj_close is stored as text in this sheet (shared string "0.680000007152557").
$ws.Range("G124").NumberFormat = "@"
$ws.Range("G124").Value = "0.680000007152557"
$ws.Range("G124").Style = "Normal"

$ws.Range("H124").Value = "BWZ.MI"

# --- Row 125 ---
$ws.Range("A123").Copy()
$ws.Range("A125").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A125").Value = 45502.2987615741

$ws.Range("B125").Value = 50
$ws.Range("C125").Value = 0.670000016689301
$ws.Range("D125").Value = 0.670000016689301
$ws.Range("E125").Value = 0.670000016689301
$ws.Range("F125").Value = 0.670000016689301

$ws.Range("G125").NumberFormat = "@"
$ws.Range("G125").Value = "0.670000016689301"
$ws.Range("G125").Style = "Normal"

$ws.Range("H125").Value = "BWZ.MI"

$excel.CutCopyMode = 0
